# "add models for silicon only"
# Adds a new "Si" (silicon-only) results section below the existing
# "CARBON, SILICON and NITROGEN" section, and re-sorts the three
# previously-sorted sections (CARBON and SILICON / CARBON and NITROGEN /
# C, Si, N) ascending by rating, matching the author's workflow of
# re-running Sort after adding the new group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Add the new "Si" section (rows 47-53), directly below the existing
#    "C, Si, N" section (which ends at row 45, with row 46 left blank as
#    a separator, matching the sheet's existing layout convention).
# ---------------------------------------------------------------------

# Section header (underlined, like the other section headers) plus the
# record-count annotation next to it.
$ws.Range("A47").Value = "Si"
$ws.Range("A47").Font.Underline = 2
$ws.Range("C47").Value = 14679
$ws.Range("C47").NumberFormat = "#,##0"
$ws.Range("D47").Value = "# records"

# Model rows, in the order they were entered (unsorted at first --
# only the older sections get re-sorted below).
$ws.Range("A48").Value = "Decision Tree"
$ws.Range("B48").Value = 80.8

$ws.Range("A49").Value = "Random Forest Classifier"
$ws.Range("B49").Value = 85.6

$ws.Range("A50").Value = "SVM radial"
$ws.Range("B50").Value = 86.8

$ws.Range("A51").Value = "Logistic Regression"
$ws.Range("B51").Value = 88.1

$ws.Range("A52").Value = "Deep Learning"
$ws.Range("B52").Value = 88.7

$ws.Range("A53").Value = "kNN"
$ws.Range("B53").Value = 89.2

# ---------------------------------------------------------------------
# 2. Re-sort the three sections (ascending by rating, column B) -- same
#    action as before, just re-run on each group.
# ---------------------------------------------------------------------

$s1 = $ws.Sort
$s1.SortFields.Clear()
$s1.SortFields.Add($ws.Range("B20:B24"))
$s1.SetRange($ws.Range("A20:B24"))
$s1.Header = 0
$s1.Apply()

$s2 = $ws.Sort
$s2.SortFields.Clear()
$s2.SortFields.Add($ws.Range("B36:B38"))
$s2.SetRange($ws.Range("A36:B38"))
$s2.Header = 0
$s2.Apply()

# This is the last sort run by the author, so it's the one that ends up
# recorded in the sheet's <sortState>.
$s3 = $ws.Sort
$s3.SortFields.Clear()
$s3.SortFields.Add($ws.Range("B41:B45"))
$s3.SetRange($ws.Range("A41:B45"))
$s3.Header = 0
$s3.Apply()

# ---------------------------------------------------------------------
# 3. Update the view: scrolled down to the new section, selection left
#    on the new last-used cell.
# ---------------------------------------------------------------------

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 38
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C53").Select()
